$d = $word.ActiveDocument

# Common run properties (Times New Roman, 12pt) used throughout this document.
$rpr = "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>"

# ---------------------------------------------------------------------------
# Edit 1: "DSC 510 Weather Forecast API" -> split into 4 runs reading
# "DSC 5" + "00" + " " + "Research Paper on Image Classification"
# (i.e. the project title becomes "DSC 500 Research Paper on Image
# Classification", but written out across 4 separate runs as in the diff).
# ---------------------------------------------------------------------------
$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r`a") -eq "DSC 510 Weather Forecast API") {
        $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
               "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>" + `
               "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" + `
               "<w:r>$rpr<w:t>DSC 5</w:t></w:r>" + `
               "<w:r>$rpr<w:t>00</w:t></w:r>" + `
               "<w:r>$rpr<w:t xml:space='preserve'> </w:t></w:r>" + `
               "<w:r>$rpr<w:t>Research Paper on Image Classification</w:t></w:r>" + `
               "</w:p>"
        $para.Range.InsertXML($xml)
        $found1 = $true
        break
    }
}
Write-Host "Edit1 (DSC 510 -> DSC 500 Research Paper):" $found1

# ---------------------------------------------------------------------------
# Edit 2: merge the runs (and drop the proofErr bookmarks) that spell out
# "There is anecdotal evidence ... to developing " into a single run. This
# edit starts at the very first run of its paragraph, so a normal
# Find/Replace over identical text naturally coalesces just those runs
# without disturbing the unrelated trailing run ("countries (i.e.
# Argentina) ...").
# ---------------------------------------------------------------------------
$old2 = "There is anecdotal evidence that males are more likely to commit suicide than females.  There were also comparisons between developed countries (i.e. United States) to developing "
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)
Write-Host "Edit2 (There is anecdotal evidence... merge runs):" $found2

# ---------------------------------------------------------------------------
# Edit 3: "3." stays its own run; the following tab + "Sex, " run absorbs
# "male" and " and female" (dropping the proofErr bookmarks around "male")
# so the paragraph ends up as two runs: "3." and <tab/>"Sex, male and
# female". A plain Find/Replace spanning those runs also pulls in the
# untouched "3." run in this runtime, so rebuild the paragraph explicitly
# via InsertXML instead.
# ---------------------------------------------------------------------------
$found3 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd("`r`a")
    if ($text -eq "3.`tSex, male and female") {
        $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
               "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='360'/><w:contextualSpacing/>" + `
               "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" + `
               "<w:r>$rpr<w:t>3.</w:t></w:r>" + `
               "<w:r>$rpr<w:tab/><w:t>Sex, male and female</w:t></w:r>" + `
               "</w:p>"
        $para.Range.InsertXML($xml)
        $found3 = $true
        break
    }
}
Write-Host "Edit3 (Sex, male and female merge runs):" $found3
